# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to "In Translation"
#   everywhere it appears (Overview!E2:F2 and the per-locale Status column C2 on the
#   "zh-cn" and "de-de" sheets).
# - Shrink the Status-related columns to match the new (shorter) status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C1").ColumnWidth = 12.5
